# Update countries & provincias Spain
# - Refresh the "last updated" timestamp
# - Update case counts for several countries (Estados Unidos, Alemania,
#   Austria, Bielorrusia, Serbia, Argentina, Islandia)
# - Re-sort the country table by "Casos totales" (column B) descending,
#   since updated totals change the ranking of a couple of countries
#   (Bielorrusia and Serbia move up in the list).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "Datos actualizados..." timestamp cell (A1) -----------
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 15:52"

# --- 2. Update per-country figures ---------------------------------------
# Map: country name -> hashtable of column letter -> new value
$updates = @{
    "Estados Unidos" = @{ "D" = 89141; "E" = 748198 }
    "Alemania"        = @{ "B" = 153544; "C" = 415; "E" = 41167; "G" = 2; "H" = 5577 }
    "Austria"         = @{ "E" = 2669; "G" = 8; "H" = 530 }
    "Bielorrusia"     = @{ "B" = 8773; "C" = 751; "D" = 1120; "E" = 7590; "F" = 92; "G" = 3; "H" = 63 }
    "Serbia"          = @{ "B" = 7483; "C" = 207; "D" = 1094; "E" = 6245; "F" = 96; "G" = 5; "H" = 144 }
    "Argentina"       = @{ "D" = 976; "E" = 2292; "G" = 2; "H" = 167 }
    "Islandia"        = @{ "D" = 1542; "E" = 237 }
}

$dataRange = $ws.Range("A4:A216")

foreach ($country in $updates.Keys) {
    # LookAt = 1 (xlWhole) so e.g. "Estados Unidos" doesn't partial-match
    # "Islas Virgenes de los Estados Unidos".
    $cell = $dataRange.Find($country, $null, $null, 1)
    $row = $cell.Row
    $cols = $updates[$country]
    foreach ($col in $cols.Keys) {
        $ws.Range($col + $row).Value = $cols[$col]
    }
}

# --- 3. Re-sort the country table by Casos totales (column B) descending -
$sortRange = $ws.Range("A4:H216")
$keyRange = $ws.Range("B4:B216")
$sortRange.Sort($keyRange, 2)
